$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Merge G19:L19 (new merged region for this logbook entry's "Logbook updated" note)
$ws.Range("G19:L19").Merge()

# Restore/apply the correct cell formatting (the native Merge() call stamps an
# outlined-border style on the merged block; copy the real formats back in
# from sibling cells that already carry the formats we want).
$ws.Range("F19").Copy()
$ws.Range("G19").PasteSpecial(-4122)
$ws.Range("M20").Copy()
$ws.Range("H19:L19").PasteSpecial(-4122)
$ws.Range("M18").Copy()
$ws.Range("M19").PasteSpecial(-4122)

# Fill in the new logbook row (row 19) values
$ws.Range("A19").Value = "Group 41"
$ws.Range("B19").Value = 43018
$ws.Range("C19").Value = 43018.020833333336
$ws.Range("D19").Value = "HV, AR, UD"
$ws.Range("E19").Value = "5 mins"
$ws.Range("F19").Value = "UP02"
$ws.Range("G19").Value = "Logbook updated"
$ws.Range("M19").Value = 0
$ws.Range("N19").Value = "0"
$ws.Range("P19").Value = "40 of 76 tests fail"
$ws.Range("Q19").Value = "UP02 with logbook updated"

# O18/O19 share the same underlying text ("29" -> "36"); update both cells to
# the new figure so the shared string itself is renamed rather than forking a
# second, duplicate string.
$ws.Range("O18").Value = "36"
$ws.Range("O19").Value = "36"
